$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the most recent day sheet (20191111) to seed the new sheet
#    with identical formulas / styles / column widths, then rename + place
#    it as the last sheet (right after 20191111).
# ---------------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("20191111")
$wsPrev.Copy([System.Reflection.Missing]::Value, $wsPrev)
$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "20191117"

# ---------------------------------------------------------------------------
# 2. Overwrite the new sheet's data (rows 1-15 position buckets, rows 20-24
#    directional summary) with the 2019-11-17 figures. Formulas (E1:E13
#    shared formula, F14, B15, C15, B24, C24) are already present from the
#    copy and simply recompute from the new inputs.
# ---------------------------------------------------------------------------
$wsNew.Range("B1").Value = 922
$wsNew.Range("C1").Value = 0.07
$wsNew.Range("D1").Value = 0

$wsNew.Range("B2").Value = 381
$wsNew.Range("C2").Value = 0.03
$wsNew.Range("D2").Value = 5

$wsNew.Range("B3").Value = 441
$wsNew.Range("C3").Value = 0.03
$wsNew.Range("D3").Value = 15

$wsNew.Range("B4").Value = 480
$wsNew.Range("C4").Value = 0.03
$wsNew.Range("D4").Value = 25

$wsNew.Range("B5").Value = 491
$wsNew.Range("C5").Value = 0.03
$wsNew.Range("D5").Value = 35

$wsNew.Range("B6").Value = 560
$wsNew.Range("C6").Value = 0.04
$wsNew.Range("D6").Value = 45

$wsNew.Range("B7").Value = 806
$wsNew.Range("C7").Value = 0.06
$wsNew.Range("D7").Value = 55

$wsNew.Range("B8").Value = 646
$wsNew.Range("C8").Value = 0.05
$wsNew.Range("D8").Value = 65

$wsNew.Range("B9").Value = 754
$wsNew.Range("C9").Value = 0.06
$wsNew.Range("D9").Value = 75

$wsNew.Range("B10").Value = 962
$wsNew.Range("C10").Value = 0.07
$wsNew.Range("D10").Value = 85

$wsNew.Range("B11").Value = 1396
$wsNew.Range("C11").Value = 0.11
$wsNew.Range("D11").Value = 95

$wsNew.Range("B12").Value = 1323
$wsNew.Range("C12").Value = 0.1
$wsNew.Range("D12").Value = 100

$wsNew.Range("B13").Value = 879
$wsNew.Range("C13").Value = 0.07
$wsNew.Range("D13").Value = 100

$wsNew.Range("B14").Value = 2242
$wsNew.Range("C14").Value = 0.18

$wsNew.Range("B20").Value = 3438
$wsNew.Range("C20").Value = 0.28

$wsNew.Range("B21").Value = 4960
$wsNew.Range("C21").Value = 0.4

$wsNew.Range("B22").Value = 1544
$wsNew.Range("C22").Value = 0.12

$wsNew.Range("B23").Value = 2335
$wsNew.Range("C23").Value = 0.19

# ---------------------------------------------------------------------------
# 3. Update the new sheet's view: scroll so row 5 is at top, and select
#    A20:C23 (active cell A20) as the last user action on this sheet.
# ---------------------------------------------------------------------------
$wsNew.Activate()
$wsNew.Range("A20:C23").Select()

# ---------------------------------------------------------------------------
# 4. The previous "current" sheet (20191111) keeps its own last selection,
#    but is no longer the active tab (handled automatically once 20191117
#    becomes active above).
# ---------------------------------------------------------------------------
$wsPrev.Range("E21").Select()

# ---------------------------------------------------------------------------
# 5. Sheet 20191110's selection moved from H34 to I34.
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("20191110")
$ws10.Activate()
$ws10.Range("I34").Select()

# ---------------------------------------------------------------------------
# 6. Re-activate the new sheet so it ends up as the active tab, matching
#    the workbook's final active-sheet state.
# ---------------------------------------------------------------------------
$wsNew.Activate()
$wsNew.Range("A20:C23").Select()
